# Daily update: append two newly-found facility name matches to the
# known_missing_matches list (fix chrome download dir bug; daily update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the two new facility names right after the current last row (A54).
$ws.Range("A55").Value = "OLNEY ES ANNEX (TRIUMPHANT FAITH-TABOR ROAD)"
$ws.Range("A56").Value = "LABRUM MS"

# Reflect the post-edit selection state (both newly added cells selected,
# active cell on the first of the two).
$ws.Range("A55:A56").Select()
